$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix value of SCALE pot: row 10 holds the SCALE1 reference, and its
# Value column (C10) incorrectly duplicated the SLOPE1 pot's value
# ("B100k"). The SCALE1 pot is actually a B10k pot.
$ws.Range("C10").Value = "B10k"

# Leave the selection where the author's edit session ended up.
$ws.Range("C11").Select()
